$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at the top of the data block (rows 2 and 3),
# pushing the existing data rows down by two.
$ws.Rows("2:3").Insert()

# The inserted rows picked up formatting from the row above (the header);
# clear that so the new data cells match the unformatted data rows below.
$ws.Range("A2:C3").ClearFormats()

# Populate the two new data rows with their values.
$ws.Range("A2").Value = 0.1914996167887811
$ws.Range("B2").Value = 0.03030422819859344
$ws.Range("C2").Value = 0.02057685541069637

$ws.Range("A3").Value = 0.1882859338884768
$ws.Range("B3").Value = 0.0655750582480559
$ws.Range("C3").Value = 0.4253946024438608
